# Remove the MAX_WAGA_PACZKI_A / _B / _C parameter rows from Parametry_Ogolne
# (rows 7-9), which shifts the former PROG_AUTOMATYCZNEGO_ZWROTU row (10) up to row 7.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Parametry_Ogolne")

$ws1.Rows.Item(7).Delete()
$ws1.Rows.Item(7).Delete()
$ws1.Rows.Item(7).Delete()

# Make Parametry_Ogolne the active sheet/tab with cell H9 selected
$ws1.Activate()
$ws1.Range("H9").Select()
